$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.327.89"
$ws.Range("E2").Value = "  -2.09%  "
$ws.Range("D3").Value = "1.564.12"
$ws.Range("E3").Value = "  -3.47%  "
$ws.Range("E4").Value = "  -0.35%  "
$ws.Range("D5").Value = "'206.90"
$ws.Range("E5").Value = "  -2.74%  "
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("D7").Value = "'0.478"
$ws.Range("E7").Value = "  -4.24%  "
$ws.Range("D8").Value = "'0.0614"
$ws.Range("E8").Value = "  -0.51%  "
$ws.Range("D9").Value = "'0.244"
$ws.Range("E9").Value = "  -2.37%  "
$ws.Range("D10").Value = "'17.89"
$ws.Range("E10").Value = "  -3.07%  "
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("D12").Value = "1.778.31"
$ws.Range("E12").Value = "  -3.60%  "
$ws.Range("D13").Value = "1.563.93"
$ws.Range("E13").Value = "  -3.46%  "
$ws.Range("E14").Value = "  -3.74%  "
$ws.Range("D15").Value = "'0.508"
$ws.Range("E15").Value = "  -2.77%  "
$ws.Range("D16").Value = "25.306.02"
$ws.Range("E16").Value = "  -2.18%  "
$ws.Range("D17").Value = "0.0₃0715"
$ws.Range("E17").Value = "  -3.05%  "
$ws.Range("D18").Value = "'59.43"
$ws.Range("E18").Value = "  -3.40%  "
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").Value = "'187.50"
$ws.Range("E20").Value = "  -2.19%  "
$ws.Range("D21").Value = "'4.13"
$ws.Range("E21").Value = "  -2.57%  "
$ws.Range("D22").Value = "'9.29"
$ws.Range("E22").Value = "  -1.98%  "
$ws.Range("D23").Value = "'5.88"
$ws.Range("E23").Value = "  -2.44%  "
$ws.Range("D24").Value = "'0.131"
$ws.Range("E24").Value = "  -3.16%  "
$ws.Range("D25").Value = "'141.03"
$ws.Range("E25").Value = "  -2.02%  "
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("E27").Value = "  -2.36%  "
$ws.Range("D28").Value = "'14.95"
$ws.Range("E28").Value = "  -1.39%  "
$ws.Range("D29").Value = "'6.40"
$ws.Range("E29").Value = "  -4.00%  "
$ws.Range("E30").Value = "  -6.80%  "
$ws.Range("D31").Value = "'0.0467"
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("D32").Value = "'3.07"
$ws.Range("E32").Value = "  -1.82%  "
$ws.Range("D33").Value = "'2.99"
$ws.Range("E33").Value = "  -3.26%  "
$ws.Range("D34").Value = "'1.49"
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("D35").Value = "'2.30"
$ws.Range("E35").Value = "  -4.25%  "
$ws.Range("D36").Value = "1.091.48"
$ws.Range("E36").Value = "  -2.86%  "
$ws.Range("B37").Value = "PaxDollar"
$ws.Range("C37").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "'2.34"
$ws.Range("E38").Value = "  -1.11%  "
$ws.Range("D39").Value = "'0.498"
$ws.Range("E39").Value = "  -2.47%  "
$ws.Range("D40").Value = "'0.0149"
$ws.Range("E40").Value = "  -2.65%  "
$ws.Range("D41").Value = "'0.776"
$ws.Range("E41").Value = "  -7.43%  "
$ws.Range("D42").Value = "'0.800"
$ws.Range("E42").Value = "  +6.80%  "
$ws.Range("D43").Value = "'93.09"
$ws.Range("D44").Value = "'5.12"
$ws.Range("E44").Value = "  +1.68%  "
$ws.Range("D45").Value = "1.693.94"
$ws.Range("E45").Value = "  -3.50%  "
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("E47").Value = "  -1.84%  "
$ws.Range("D48").Value = "'52.62"
$ws.Range("E48").Value = "  -2.62%  "
$ws.Range("D49").Value = "'0.0505"
$ws.Range("E49").Value = "  -2.98%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.404"
$ws.Range("E50").Value = "  -1.86%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  -0.38%  "
